$wb = $excel.ActiveWorkbook

# The original sheet ("Sayfa1") loses its header row content.
$sheet1 = $wb.Worksheets.Item(1)
$sheet1.Cells.Clear()
$sheet1.Range("H8").Select()

# Add a new worksheet named "Stock" positioned right after "Sayfa1"
$stock = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$stock.Name = "Stock"

# Match the sheetPr block (outlinePr) of the original sheet
$stock.Outline.SummaryColumn = 1
$stock.Outline.SummaryRow = 1

# Page margins
$stock.PageSetup.LeftMargin = 54
$stock.PageSetup.RightMargin = 54
$stock.PageSetup.TopMargin = 72
$stock.PageSetup.BottomMargin = 72
$stock.PageSetup.HeaderMargin = 36
$stock.PageSetup.FooterMargin = 36

# Header row
$stock.Range("A1").Value = "TÜR"
$stock.Range("B1").Value = "ADET"

# Data rows (quantities stored as text, matching source data)
$stock.Range("B2:B3").NumberFormat = "@"

$stock.Range("A2").Value = "muz"
$stock.Range("B2").Value = "12"

$stock.Range("A3").Value = "elma"
$stock.Range("B3").Value = "30"

# Drop the temporary "Text" number format so the cells keep default styling
$stock.Range("B2:B3").Style = "Normal"

$stock.Range("A1").Select()

# Keep "Sayfa1" as the active/selected sheet
$sheet1.Activate()
